$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as TEXT (mirrors the workbook's original inlineStr cells).
# A leading apostrophe forces Excel to store the string verbatim instead of
# renormalizing numeric-looking text (e.g. "1.00" -> 1, "18.30" -> 18.3).
function Set-TextValue($row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

Set-TextValue 2 4 '70.843.57'
Set-TextValue 2 5 '  +1.77%  '
Set-TextValue 3 4 '3.632.63'
Set-TextValue 3 5 '  +4.06%  '
Set-TextValue 4 4 '0.999'
Set-TextValue 4 5 '  +0.11%  '
Set-TextValue 5 4 '605.01'
Set-TextValue 5 5 '  +0.30%  '
Set-TextValue 6 4 '200.59'
Set-TextValue 6 5 '  +3.47%  '
Set-TextValue 7 4 '0.628'
Set-TextValue 7 5 '  +0.46%  '
Set-TextValue 8 5 '  +0.08%  '
Set-TextValue 9 4 '0.218'
Set-TextValue 9 5 '  +9.31%  '
Set-TextValue 10 5 '  -0.27%  '
Set-TextValue 11 4 '53.86'
Set-TextValue 11 5 '  +1.33%  '
Set-TextValue 12 5 '  +2.29%  '
Set-TextValue 13 4 '9.57'
Set-TextValue 13 5 '  +1.35%  '
Set-TextValue 14 4 '4.206.71'
Set-TextValue 14 5 '  +3.44%  '
Set-TextValue 15 4 '629.24'
Set-TextValue 15 5 '  +6.07%  '
Set-TextValue 16 5 '  +2.10%  '
Set-TextValue 17 4 '70.885.27'
Set-TextValue 17 5 '  +1.62%  '
Set-TextValue 18 4 '3.609.73'
Set-TextValue 18 5 '  +3.18%  '
Set-TextValue 19 5 '  +0.64%  '
Set-TextValue 20 5 '  +0.98%  '
Set-TextValue 21 4 '1.00'
Set-TextValue 21 5 '  +1.73%  '
Set-TextValue 22 4 '18.30'
Set-TextValue 22 5 '  +1.39%  '
Set-TextValue 23 4 '5.33'
Set-TextValue 23 5 '  +0.38%  '
Set-TextValue 24 4 '104.02'
Set-TextValue 24 5 '  +2.07%  '
Set-TextValue 25 5 '  -0.47%  '
Set-TextValue 26 4 '3.01'
Set-TextValue 26 5 '  -4.19%  '
Set-TextValue 27 5 '  -2.20%  '
Set-TextValue 28 5 '  +2.54%  '
Set-TextValue 29 5 '  +1.52%  '
Set-TextValue 30 4 '4.75'
Set-TextValue 30 5 '  +13.89%  '
Set-TextValue 31 5 '  +2.89%  '
Set-TextValue 32 4 '12.25'
Set-TextValue 32 5 '  -0.90%  '
Set-TextValue 33 5 '  +1.63%  '
Set-TextValue 34 4 '63.37'
Set-TextValue 34 5 '  +0.41%  '
Set-TextValue 35 4 '0.0₃0881'
Set-TextValue 35 5 '  +6.17%  '
Set-TextValue 36 4 '3.990.06'
Set-TextValue 36 5 '  +7.58%  '
Set-TextValue 37 5 '  +0.17%  '
Set-TextValue 38 4 '515.11'
Set-TextValue 38 5 '  +7.63%  '
Set-TextValue 39 4 '3.03'
Set-TextValue 39 5 '  -1.79%  '
Set-TextValue 40 2 'InjectiveProtocol'
Set-TextValue 40 3 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 40 4 '36.77'
Set-TextValue 40 5 '  +1.27%  '
Set-TextValue 41 2 'TheGraph'
Set-TextValue 41 3 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 41 4 '0.390'
Set-TextValue 41 5 '  +0.20%  '
Set-TextValue 42 5 '  -2.58%  '
Set-TextValue 43 5 '  +3.32%  '
Set-TextValue 44 5 '  +2.29%  '
Set-TextValue 45 5 '  +6.61%  '
Set-TextValue 46 4 '2.91'
Set-TextValue 46 5 '  +4.16%  '
Set-TextValue 47 5 '  +1.00%  '
Set-TextValue 48 4 '8.61'
Set-TextValue 48 5 '  +2.48%  '
Set-TextValue 49 5 '  -0.26%  '
Set-TextValue 50 4 '0.000251'
Set-TextValue 50 5 '  +2.99%  '
Set-TextValue 51 5 '  +1.76%  '
